# Weekly update: insert the latest week's row of data at the top of the
# Coliflor (Terminal Hortofrutícola Agro Chillán) price series (row 68,
# just after the header/fixed leading rows), pushing the existing history
# down by one row. The workbook keeps a rolling weekly series, so last
# week's brand-new row becomes this week's newest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 68:175 down to 69:176, carrying values/formats along.
$ws.Rows("68:68").Insert()

# Populate the new row 68 with this week's record.
$ws.Cells.Item(68, 1).Value  = 7
$ws.Cells.Item(68, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(68, 3).Value  = "Ñuble"
$ws.Cells.Item(68, 4).Value  = 44495
$ws.Cells.Item(68, 5).Value  = 16
$ws.Cells.Item(68, 6).Value  = 100112008
$ws.Cells.Item(68, 7).Value  = "Coliflor"
$ws.Cells.Item(68, 8).Value  = "Sin especificar"
$ws.Cells.Item(68, 9).Value  = "Primera"
$ws.Cells.Item(68, 10).Value = 160
$ws.Cells.Item(68, 11).Value = 700
$ws.Cells.Item(68, 12).Value = 800
$ws.Cells.Item(68, 13).Value = 750
$ws.Cells.Item(68, 14).Value = "`$/unidad"
$ws.Cells.Item(68, 15).Value = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value = 750
$ws.Cells.Item(68, 17).Value = 1
$ws.Cells.Item(68, 18).Value = "Hortaliza"

# Note: the inserted row already inherits the date number-format style used
# by the rest of column D from the row below, so no extra style copy needed.
